# Apply the commit's edits to the "rotas" worksheet:
#  - update delivery dates (column B) for rows 2-20
#  - fix a restriction time (O2)
#  - fix CEP / house-number typos (C5, C8, D8)
#
# NOTE: several of these cells hold dates as plain TEXT (t="inlineStr"),
# not real date serials. Assigning a date-looking string straight to
# Range.Value/.Value2/.Formula makes Excel "helpfully" reinterpret it as
# a date serial (and stamp a date NumberFormat on the cell), which would
# change the cell's stored type/style and not match the source data.
# To keep these as literal text, we write the text as a quoted formula
# (=\"11/4/2024\") so it is never date-parsed, then convert that formula
# to a plain static value in place via Copy / PasteSpecial values-only,
# which leaves no formula behind and does not touch the cell's style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    $escaped = $Text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy($cell) | Out-Null
    $cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues) | Out-Null
}

$excel.CutCopyMode = $false

# Column B - Data Entrega (delivery date), stored as text "M/D/YYYY"
Set-TextValue "B2"  "11/4/2024"
Set-TextValue "B3"  "11/4/2024"
Set-TextValue "B4"  "11/4/2024"
Set-TextValue "B5"  "11/4/2024"
Set-TextValue "B6"  "11/4/2024"
Set-TextValue "B7"  "11/5/2024"
Set-TextValue "B8"  "11/5/2024"
Set-TextValue "B9"  "11/5/2024"
Set-TextValue "B10" "11/5/2024"
Set-TextValue "B11" "11/5/2024"
Set-TextValue "B12" "11/6/2024"
Set-TextValue "B13" "11/6/2024"
Set-TextValue "B14" "11/6/2024"
Set-TextValue "B15" "11/6/2024"
Set-TextValue "B16" "11/5/2024"
Set-TextValue "B17" "11/5/2024"
Set-TextValue "B18" "11/4/2024"
Set-TextValue "B19" "11/6/2024"
Set-TextValue "B20" "11/6/2024"

# O2 - restriction time text "8:00" -> "16:00"
$ws.Range("O2").Value = "16:00"

# Fix typo'd CEP / house number values
$ws.Range("C5").Value = 88070142
$ws.Range("C8").Value = 88106550
$ws.Range("D8").Value = 150

$excel.CutCopyMode = $false
